# Apply the "Removing past portfolio data" update:
#  - Shift every timestamp in column A (rows 2-97) forward by 4 days
#    (the data window moved from 2025-03-20 to 2025-03-24)
#  - Replace the production values in column B (rows 2-97) with the
#    refreshed values pulled from the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 97

# Shift all timestamps in column A by +4 days, preserving their existing
# formatting/style.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldSerial = [double]$cell.Value2
    $cell.Value = $oldSerial + 4
}

# New column B values (Actual Production, MW) for rows 2..97.
$bValues = @(
    1833,1818,1775,1650,1571,1494,1381,1279,1149,1024,907,849,867,837,746,
    732,789,826,871,868,839,845,828,799,716,674,614,552,467,448,428,385,
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
}
